$wb = $excel.ActiveWorkbook

# The existing "SourceData" sheet originally held a generic value-coding
# schema (code/display/value_code/...). We are introducing a new
# "DataSource" sheet (connection details for a data source) that sits right
# before "SourceData" in the tab order, and reworking "SourceData" itself
# into a lean reference sheet that points at a DataSource.

$sourceDataRef = $wb.Worksheets.Item("SourceData")

# Insert the new sheet immediately before "SourceData".
$dataSource = $wb.Worksheets.Add($sourceDataRef)
$dataSource.Name = "DataSource"

# Match the page-setup margins used throughout the rest of the workbook.
$dataSource.PageSetup.LeftMargin = 54
$dataSource.PageSetup.RightMargin = 54
$dataSource.PageSetup.TopMargin = 72
$dataSource.PageSetup.BottomMargin = 72
$dataSource.PageSetup.HeaderMargin = 36
$dataSource.PageSetup.FooterMargin = 36

$dataSource.Range("A1").Value = "snapshot_id"
$dataSource.Range("B1").Value = "google_data_project"
$dataSource.Range("C1").Value = "snapshot_dataset"
$dataSource.Range("D1").Value = "table"
$dataSource.Range("E1").Value = "parameterized_query"
$dataSource.Range("F1").Value = "id"
$dataSource.Range("G1").Value = "external_id"

# Re-fetch the "SourceData" worksheet by name now that a new sheet has been
# inserted in front of it (the old handle can end up referring to the newly
# inserted sheet instead).
$sourceData = $wb.Worksheets.Item("SourceData")

# Replace the old SourceData header row with the new, smaller schema.
$sourceData.Cells.Clear()

$sourceData.Range("A1").Value = "data_source"
$sourceData.Range("B1").Value = "query_parameter"
$sourceData.Range("C1").Value = "has_access_policy"
$sourceData.Range("D1").Value = "id"
$sourceData.Range("E1").Value = "external_id"

# Restore the original active sheet (the new sheet becomes active on
# insertion, which would otherwise shift the workbook's saved tab
# selection away from "Subject").
$wb.Worksheets.Item("Subject").Activate()
